$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing the existing agenda rows
# (old rows 10-14) down to rows 11-15.
$ws.Rows.Item(10).Insert()

# The freshly-inserted row has no formatting of its own; copy the
# look of the row above (borders/font) onto it so it matches the
# rest of the agenda table, then overwrite the inherited text.
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row 10 with the new agenda item (Speaker | Subject | Time)
$ws.Range("A10").Value = "Sammy Douglas"
$ws.Range("B10").Value = "Destruction Update"
$ws.Range("C10").Value = "5 minutes"

# Update the selection to match the saved workbook state
$ws.Range("C11").Select()
